# Swap the deck's colour theme: the design ("Integral") that is currently
# applied to the slide master / presentation (backed by ppt/theme/theme2.xml)
# is restored to the stock "Office Theme" palette, i.e. the 12 theme colour
# slots (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) get the default
# Office colour values instead of the Integral ones.
#
# PowerPoint keeps only one ThemeColorScheme per slide master, reachable off
# any Slide/SlideRange/NotesPage object, and editing its RGB values edits the
# theme part (ppt/theme/theme2.xml) backing the slide master/presentation in
# place - exactly the part that needs to change.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$cs = $s.ThemeColorScheme

# Office Theme colours (standard Office default palette), in VBA RGB()
# encoding (R + G*256 + B*65536), ordered to match ThemeColorScheme.Item(1..12):
#   1 dk1=000000  2 lt1=FFFFFF  3 dk2=44546A  4 lt2=E7E6E6
#   5 accent1=5B9BD5  6 accent2=ED7D31  7 accent3=A5A5A5  8 accent4=FFC000
#   9 accent5=4472C4  10 accent6=70AD47  11 hlink=0563C1  12 folHlink=954F72
$officeThemeColors = @(
    0,        # dk1      000000
    16777215, # lt1      FFFFFF
    6968388,  # dk2      44546A
    15132391, # lt2      E7E6E6
    13998939, # accent1  5B9BD5
    3243501,  # accent2  ED7D31
    10855845, # accent3  A5A5A5
    49407,    # accent4  FFC000
    12874308, # accent5  4472C4
    4697456,  # accent6  70AD47
    12673797, # hlink    0563C1
    7491477   # folHlink 954F72
)

for ($i = 1; $i -le 12; $i++) {
    $cs.Item($i).RGB = $officeThemeColors[$i - 1]
}
